$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
